$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 28 (J3/DAC_CLK_P row): pin number corrected from 10 to text "c"
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = "c"

# ---------------------------------------------------------------------------
# New rows describing extra J3/J4 pins (DFT / sigma-delta ADC additions).
# The text cells below are written in the exact order the strings were first
# introduced so the shared-string table order matches the source edit.
# ---------------------------------------------------------------------------
$ws.Range("E35").Value = "~CW"
$ws.Range("D38").Value = "h14"
$ws.Range("E37").Value = "CLK_EXT_IN"
$ws.Range("D35").Value = "m15"
$ws.Range("D37").Value = "h16"
$ws.Range("E38").Value = "CLK_OUT"
$ws.Range("E40").Value = "LED0"
$ws.Range("E41").Value = "LED1"
$ws.Range("E42").Value = "LED2"
$ws.Range("E43").Value = "LED3"
$ws.Range("D40").Value = "f9"
$ws.Range("D41").Value = "e8"
$ws.Range("D42").Value = "e7"
$ws.Range("D43").Value = "d7"

# ---------------------------------------------------------------------------
# Remaining (non shared-string-table-affecting) cell values: connector name
# and pin number columns.
# ---------------------------------------------------------------------------
$ws.Range("B35").Value = "J4"
$ws.Range("C35").Value = 9

$ws.Range("B37").Value = "J3"
$ws.Range("C37").Value = 25

$ws.Range("B38").Value = "J3"
$ws.Range("C38").Value = 19

$ws.Range("B40").Value = "J3"
$ws.Range("C40").Value = 13

$ws.Range("B41").Value = "J3"
$ws.Range("C41").Value = 15

$ws.Range("B42").Value = "J3"
$ws.Range("C42").Value = 17

$ws.Range("B43").Value = "J3"
$ws.Range("C43").Value = 19

# ---------------------------------------------------------------------------
# Update selection / view to match the final saved state
# ---------------------------------------------------------------------------
$ws.Range("D43").Select()
